# Updates scraped_at timestamps on the "snapshot" sheet (column K, rows 2-39)
# and removes the newly-injured player row that has since been processed on
# the "new_injured" sheet (row 2), matching a re-run of the KHL injuries scraper.

$wb = $excel.ActiveWorkbook

# --- 1. "snapshot" sheet: refresh the scraped_at (column K) timestamps ---
$snapshot = $wb.Worksheets.Item("snapshot")

$newScrapedAt = @(
    "2025-12-10T07:01:17.183136+00:00",
    "2025-12-10T07:01:17.183168+00:00",
    "2025-12-10T07:01:17.183188+00:00",
    "2025-12-10T07:01:19.953163+00:00",
    "2025-12-10T07:01:19.953191+00:00",
    "2025-12-10T07:01:19.953208+00:00",
    "2025-12-10T07:01:22.280323+00:00",
    "2025-12-10T07:01:25.032611+00:00",
    "2025-12-10T07:01:27.894934+00:00",
    "2025-12-10T07:01:30.631402+00:00",
    "2025-12-10T07:01:35.848327+00:00",
    "2025-12-10T07:01:35.848370+00:00",
    "2025-12-10T07:01:38.560321+00:00",
    "2025-12-10T07:01:41.388169+00:00",
    "2025-12-10T07:01:44.180151+00:00",
    "2025-12-10T07:01:46.851051+00:00",
    "2025-12-10T07:01:46.851082+00:00",
    "2025-12-10T07:01:49.239314+00:00",
    "2025-12-10T07:01:49.239343+00:00",
    "2025-12-10T07:01:49.239360+00:00",
    "2025-12-10T07:01:51.912124+00:00",
    "2025-12-10T07:01:51.912154+00:00",
    "2025-12-10T07:01:51.912174+00:00",
    "2025-12-10T07:01:51.912192+00:00",
    "2025-12-10T07:01:51.912209+00:00",
    "2025-12-10T07:01:54.725460+00:00",
    "2025-12-10T07:01:54.725496+00:00",
    "2025-12-10T07:01:54.725518+00:00",
    "2025-12-10T07:01:56.995991+00:00",
    "2025-12-10T07:02:12.499362+00:00",
    "2025-12-10T07:02:12.499390+00:00",
    "2025-12-10T07:02:12.499408+00:00",
    "2025-12-10T07:02:14.760136+00:00",
    "2025-12-10T07:02:14.760165+00:00",
    "2025-12-10T07:02:17.634138+00:00",
    "2025-12-10T07:02:17.634167+00:00",
    "2025-12-10T07:02:20.376394+00:00",
    "2025-12-10T07:02:20.376422+00:00"
)

for ($i = 0; $i -lt $newScrapedAt.Length; $i++) {
    $row = $i + 2
    $snapshot.Cells.Item($row, 11).Value = $newScrapedAt[$i]
}

# --- 2. "new_injured" sheet: drop the row that is no longer newly-injured ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()

Write-Host "Updated $($newScrapedAt.Length) scraped_at values and removed 1 row from new_injured"
